$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 13 values (B13, C13, D13) ---
$ws.Range("B13").Value = -0.040456817157672728
$ws.Range("C13").Value = 0.11504492669835839
$ws.Range("D13").Value = 0.057729375378128042

# --- Row 14: give it the same "separator" look as rows 1-13 (height/format),
#     then update B14's value. Copy the format from row 13 (formats only) so
#     the cell style indices (date / percent styles) are reused rather than
#     duplicated, then overwrite the values back to what row 14 should hold. ---
$ws.Range("A13:D13").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Range("A14").Value = 45670
$ws.Range("B14").Value = 0.081189267611627647
$ws.Range("C14").Value = 0.01217049583297181
$ws.Range("D14").Value = -0.047139261276747425

# --- Row 15: used to hold just a lone formatted D15 cell; now gets a full
#     data row. Reuse row 14's per-column formats (date / percent) via
#     copy+paste-formats so style indices match, then set the new values. ---
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 45701

$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = -0.20628550570397458

$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = -0.14317981410163677

$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = -0.041385216162853962

# --- Row 16: the lone empty formatted cell moves from C16 to D16. ---
$ws.Range("C16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("C16").Clear()

# --- Selection moves from E15 to L15. ---
$ws.Range("L15").Select()
